$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.311.96'
$ws.Range("E2").Value = '  -1.77%  '
$ws.Range("D3").Value = '1.626.37'
$ws.Range("E3").Value = '  -1.92%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").Value = "'298.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.57%  '
$ws.Range("D7").Value = "'0.3761"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.45%  '
$ws.Range("D8").Value = "'49.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.38%  '
$ws.Range("D9").Value = "'0.3470"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.98%  '
$ws.Range("D10").Value = "'0.08028"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.07%  '
$ws.Range("E11").Value = '  -3.00%  '
$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.03%  '
$ws.Range("D13").Value = "'21.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.23%  '
$ws.Range("D14").Value = "'6.265"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.96%  '
$ws.Range("D15").Value = "'7.194"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.31%  '
$ws.Range("D16").Value = "'0.00001187"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.28%  '
$ws.Range("D17").Value = '1.627.38'
$ws.Range("E17").Value = '  -2.01%  '
$ws.Range("D18").Value = "'94.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.26%  '
$ws.Range("D19").Value = "'0.06939"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.23%  '
$ws.Range("D20").Value = "'6.594"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.47%  '
$ws.Range("D21").Value = "'17.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.45%  '
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").Value = "'12.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.62%  '
$ws.Range("D24").Value = '23.326.44'
$ws.Range("E24").Value = '  -1.69%  '
$ws.Range("D25").Value = "'2.413"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.92%  '
$ws.Range("E26").Value = '  -1.63%  '
$ws.Range("D27").Value = "'20.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.22%  '
$ws.Range("D28").Value = "'149.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.78%  '
$ws.Range("D29").Value = "'5.165"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.15%  '
$ws.Range("D30").Value = "'130.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.82%  '
$ws.Range("D31").Value = '1.802.46'
$ws.Range("E31").Value = '  -2.14%  '
$ws.Range("D32").Value = "'6.697"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.12%  '
$ws.Range("D33").Value = "'2.114"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.35%  '
$ws.Range("D34").Value = "'11.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.09%  '
$ws.Range("D35").Value = "'0.9821"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.96%  '
$ws.Range("D36").Value = "'0.02646"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.11%  '
$ws.Range("D37").Value = "'0.08700"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.27%  '
$ws.Range("D38").Value = "'0.2405"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.65%  '
$ws.Range("D39").Value = "'5.781"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.51%  '
$ws.Range("D40").Value = "'0.06739"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.72%  '
$ws.Range("E41").Value = '  -2.77%  '
$ws.Range("D42").Value = "'0.6774"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.26%  '
$ws.Range("E43").Value = '  -3.86%  '
$ws.Range("D44").Value = "'15.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.84%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").Value = "'0.6276"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.88%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D47").Value = "'3.887"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.02%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = "'2.219"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.82%  '
$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").Value = "'125.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.94%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = "'0.07610"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.84%  '
$ws.Range("D51").Value = "'1.219"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.05%  '
